$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.524.83'
$ws.Range("E2").Value = '  +0.05%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.463.57'
$ws.Range("E3").Value = '  -0.87%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.41%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.29'
$ws.Range("E5").Value = '  -0.52%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '91.30'
$ws.Range("E6").Value = '  -2.41%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.549'
$ws.Range("E7").Value = '  +0.23%  '
$ws.Range("E8").Value = '  -0.42%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.512'
$ws.Range("E9").Value = '  +3.26%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.46'
$ws.Range("E10").Value = '  -2.76%  '
$ws.Range("E11").Value = '  +1.16%  '
$ws.Range("E12").Value = '  +0.45%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.849.44'
$ws.Range("E13").Value = '  -0.93%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.83'
$ws.Range("E14").Value = '  -1.01%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.81'
$ws.Range("E15").Value = '  +2.54%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.456.56'
$ws.Range("E16").Value = '  -0.23%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.774'
$ws.Range("E17").Value = '  -1.85%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '41.510.70'
$ws.Range("E18").Value = '  +0.23%  '
$ws.Range("E19").Value = '  +2.79%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0940'
$ws.Range("E20").Value = '  +1.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '70.75'
$ws.Range("E21").Value = '  +0.78%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.06'
$ws.Range("E22").Value = '  -0.98%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.61'
$ws.Range("E23").Value = '  +0.95%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.71'
$ws.Range("E24").Value = '  -1.54%  '
$ws.Range("E25").Value = '  +0.52%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.51'
$ws.Range("E27").Value = '  +1.35%  '
$ws.Range("E28").Value = '  -0.58%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.66'
$ws.Range("E29").Value = '  -1.69%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.29'
$ws.Range("E30").Value = '  -4.03%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '155.59'
$ws.Range("E31").Value = '  +1.43%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.42'
$ws.Range("E32").Value = '  -0.99%  '
$ws.Range("E33").Value = '  +0.16%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0756'
$ws.Range("E34").Value = '  +0.24%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '17.11'
$ws.Range("E35").Value = '  -4.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.35'
$ws.Range("E36").Value = '  -6.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.87'
$ws.Range("E37").Value = '  -5.86%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.103'
$ws.Range("E38").Value = '  +3.08%  '
$ws.Range("E39").Value = '  +0.61%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.78'
$ws.Range("E40").Value = '  -4.56%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.97'
$ws.Range("E41").Value = '  -2.89%  '
$ws.Range("E42").Value = '  -0.67%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.940.01'
$ws.Range("E43").Value = '  -2.26%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0282'
$ws.Range("E44").Value = '  -0.47%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.68'
$ws.Range("E45").Value = '  -5.50%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.88'
$ws.Range("E46").Value = '  -3.41%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.03'
$ws.Range("E47").Value = '  +2.49%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.708.08'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '96.68'
$ws.Range("E49").Value = '  +0.20%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '66.57'
$ws.Range("E50").Value = '  -3.29%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '52.28'
$ws.Range("E51").Value = '  +3.54%  '
